$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4039.525
$ws.Range("I15").Value = 4039.525
$ws.Range("K15").Value = 12118.575
$ws.Range("M15").Value = -11949.575

$ws.Range("H107").Value = 334.95
$ws.Range("I107").Value = 205.21053
$ws.Range("K107").Value = 205.21053
$ws.Range("M107").Value = 1714.78947

$ws.Range("H132").Value = 1034.1094
$ws.Range("I132").Value = 934.0517
$ws.Range("K132").Value = 2802.1551
$ws.Range("M132").Value = -272.1550999999999

$ws.Range("H135").Value = 35714748
$ws.Range("I135").Value = 501.12
$ws.Range("J135").Value = 333333470
$ws.Range("K135").Value = 4510.08
$ws.Range("L135").Value = 3000001230
$ws.Range("M135").Value = -1975.08
$ws.Range("N135").Value = -3000006300

$ws.Range("H138").Value = 1547.6582
$ws.Range("I138").Value = 1217.8226
$ws.Range("J138").Value = 2750.5881
$ws.Range("K138").Value = 3653.4678
$ws.Range("L138").Value = 8251.764299999999
$ws.Range("M138").Value = 1486.5322
$ws.Range("N138").Value = -18531.7643

$ws.Range("H141").Value = 1001658.56
$ws.Range("I141").Value = 1218210.5
$ws.Range("J141").Value = 5519.6
$ws.Range("K141").Value = 3654631.5
$ws.Range("L141").Value = 16558.8
$ws.Range("M141").Value = -3649451.5
$ws.Range("N141").Value = -26918.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1111849.8
$ws.Range("I2").Value = 1389312.2
$ws.Range("K2").Value = 1389312.2
$ws.Range("M2").Value = -1389199.2

$ws.Range("H32").Value = 3222.0706
$ws.Range("I32").Value = 2588.4666
$ws.Range("K32").Value = 2588.4666
$ws.Range("M32").Value = -2301.4666

$ws.Range("H61").Value = 1316.9714
$ws.Range("I61").Value = 713.4400000000001
$ws.Range("J61").Value = 2825.8
$ws.Range("K61").Value = 713.4400000000001
$ws.Range("L61").Value = 2825.8
$ws.Range("M61").Value = -501.4400000000001
$ws.Range("N61").Value = -3249.8

$ws.Range("H74").Value = 1138.2653
$ws.Range("I74").Value = 870.19446
$ws.Range("K74").Value = 870.19446
$ws.Range("M74").Value = 3.805539999999951

$ws.Range("H77").Value = 1138.2653
$ws.Range("I77").Value = 870.19446
$ws.Range("K77").Value = 4350.9723
$ws.Range("M77").Value = 17.02769999999964

$ws.Range("H97").Value = 1036.3334
$ws.Range("I97").Value = 896.8570999999999
$ws.Range("K97").Value = 896.8570999999999
$ws.Range("M97").Value = -400.8570999999999

$ws.Range("H110").Value = 1250.7142
$ws.Range("I110").Value = 995.43335
$ws.Range("J110").Value = 2782.4
$ws.Range("K110").Value = 995.43335
$ws.Range("L110").Value = 2782.4
$ws.Range("M110").Value = 1049.56665
$ws.Range("N110").Value = -6872.4

$ws.Range("H116").Value = 1111849.8
$ws.Range("I116").Value = 1389312.2
$ws.Range("K116").Value = 1389312.2
$ws.Range("M116").Value = -1387018.2

$ws.Range("H122").Value = 1878.68
$ws.Range("I122").Value = 1493
$ws.Range("K122").Value = 4479
$ws.Range("M122").Value = -2029

$ws.Range("H132").Value = 1241.4807
$ws.Range("I132").Value = 1049.5366
$ws.Range("J132").Value = 1956.909
$ws.Range("K132").Value = 3148.6098
$ws.Range("L132").Value = 5870.727000000001
$ws.Range("M132").Value = -618.6097999999997
$ws.Range("N132").Value = -10930.727

$ws.Range("H136").Value = 1316.9714
$ws.Range("I136").Value = 713.4400000000001
$ws.Range("J136").Value = 2825.8
$ws.Range("K136").Value = 2140.32
$ws.Range("L136").Value = 8477.400000000001
$ws.Range("M136").Value = 409.6799999999998
$ws.Range("N136").Value = -13577.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1111849.8
$ws.Range("I3").Value = 1389312.2
$ws.Range("K3").Value = 1389312.2
$ws.Range("M3").Value = -1389198.2

$ws.Range("H20").Value = 2297.8333
$ws.Range("I20").Value = 2134.2727
$ws.Range("J20").Value = 2554.8572
$ws.Range("K20").Value = 2134.2727
$ws.Range("L20").Value = 2554.8572
$ws.Range("M20").Value = -1887.2727
$ws.Range("N20").Value = -3048.8572

$ws.Range("H86").Value = 431539.66
$ws.Range("I86").Value = 628363.1
$ws.Range("J86").Value = 169108.33
$ws.Range("K86").Value = 628363.1
$ws.Range("L86").Value = 169108.33
$ws.Range("M86").Value = -627240.1
$ws.Range("N86").Value = -171354.33

$ws.Range("H89").Value = 431539.66
$ws.Range("I89").Value = 628363.1
$ws.Range("J89").Value = 169108.33
$ws.Range("K89").Value = 3141815.5
$ws.Range("L89").Value = 845541.6499999999
$ws.Range("M89").Value = -3136199.5
$ws.Range("N89").Value = -856773.6499999999

$ws.Range("H107").Value = 1690.7273
$ws.Range("I107").Value = 1764.3334
$ws.Range("J107").Value = 1602.4
$ws.Range("K107").Value = 1764.3334
$ws.Range("L107").Value = 1602.4
$ws.Range("M107").Value = 155.6666
$ws.Range("N107").Value = -5442.4

$ws.Range("H134").Value = 12947.875
$ws.Range("I134").Value = 12705.077
$ws.Range("K134").Value = 38115.231
$ws.Range("M134").Value = -35580.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1683.4849
$ws.Range("I31").Value = 1265.591
$ws.Range("J31").Value = 2519.2727
$ws.Range("K31").Value = 1265.591
$ws.Range("L31").Value = 2519.2727
$ws.Range("M31").Value = -970.5909999999999
$ws.Range("N31").Value = -3109.2727

$ws.Range("H34").Value = 1683.4849
$ws.Range("I34").Value = 1265.591
$ws.Range("J34").Value = 2519.2727
$ws.Range("K34").Value = 1265.591
$ws.Range("L34").Value = 2519.2727
$ws.Range("M34").Value = -1063.591
$ws.Range("N34").Value = -2923.2727

$ws.Range("H58").Value = 1978426
$ws.Range("I58").Value = 2899445.2
$ws.Range("J58").Value = 4813.5713
$ws.Range("K58").Value = 2899445.2
$ws.Range("L58").Value = 4813.5713
$ws.Range("M58").Value = -2899242.2
$ws.Range("N58").Value = -5219.5713

$ws.Range("H107").Value = 477.14285
$ws.Range("I107").Value = 489.5
$ws.Range("K107").Value = 489.5
$ws.Range("M107").Value = 1430.5

$ws.Range("H132").Value = 1258.8959
$ws.Range("I132").Value = 812.7838
$ws.Range("K132").Value = 2438.3514
$ws.Range("M132").Value = 91.64859999999999

$ws.Range("H134").Value = 1912.641
$ws.Range("I134").Value = 1938.6
$ws.Range("J134").Value = 1866.2858
$ws.Range("K134").Value = 5815.799999999999
$ws.Range("L134").Value = 5598.857400000001
$ws.Range("M134").Value = -3280.799999999999
$ws.Range("N134").Value = -10668.8574

$ws.Range("H136").Value = 1978426
$ws.Range("I136").Value = 2899445.2
$ws.Range("J136").Value = 4813.5713
$ws.Range("K136").Value = 8698335.600000001
$ws.Range("L136").Value = 14440.7139
$ws.Range("M136").Value = -8695785.600000001
$ws.Range("N136").Value = -19540.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 194.08333
$ws.Range("I12").Value = 75
$ws.Range("J12").Value = 233.77777
$ws.Range("K12").Value = 225
$ws.Range("L12").Value = 701.33331
$ws.Range("M12").Value = -52
$ws.Range("N12").Value = -1047.33331

$ws.Range("H114").Value = 1883.4546
$ws.Range("J114").Value = 2684.2856
$ws.Range("L114").Value = 8052.8568
$ws.Range("N114").Value = -14560.8568

$ws.Range("H128").Value = 396666.34
$ws.Range("I128").Value = 396666.34
$ws.Range("K128").Value = 1189999.02
$ws.Range("M128").Value = -1185019.02

$ws.Range("H131").Value = 7948418.5
$ws.Range("I131").Value = 31250364
$ws.Range("J131").Value = 15841.574
$ws.Range("K131").Value = 93751092
$ws.Range("L131").Value = 47524.722
$ws.Range("M131").Value = -93746052
$ws.Range("N131").Value = -57604.722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3512448
$ws.Range("I11").Value = 4107103.5
$ws.Range("K11").Value = 4107103.5
$ws.Range("M11").Value = -4106964.5

$ws.Range("H70").Value = 4273.375
$ws.Range("I70").Value = 4230
$ws.Range("J70").Value = 4299.4
$ws.Range("K70").Value = 4230
$ws.Range("L70").Value = 4299.4
$ws.Range("M70").Value = -3960
$ws.Range("N70").Value = -4839.4

$ws.Range("H73").Value = 4273.375
$ws.Range("I73").Value = 4230
$ws.Range("J73").Value = 4299.4
$ws.Range("K73").Value = 4230
$ws.Range("L73").Value = 4299.4
$ws.Range("M73").Value = -3294
$ws.Range("N73").Value = -6171.4

$ws.Range("H107").Value = 262.375
$ws.Range("I107").Value = 279.8
$ws.Range("K107").Value = 279.8
$ws.Range("M107").Value = 1640.2

$ws.Range("H113").Value = 1305.0834
$ws.Range("I113").Value = 1332.625
$ws.Range("K113").Value = 1332.625
$ws.Range("M113").Value = 837.375

$ws.Range("H132").Value = 1284410.6
$ws.Range("I132").Value = 1749965.9
$ws.Range("K132").Value = 5249897.699999999
$ws.Range("M132").Value = -5247367.699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1265.9818
$ws.Range("I132").Value = 1036.4445
$ws.Range("J132").Value = 1700.8948
$ws.Range("K132").Value = 3109.3335
$ws.Range("L132").Value = 5102.6844
$ws.Range("M132").Value = -579.3335000000002
$ws.Range("N132").Value = -10162.6844

$ws.Range("H136").Value = 2327.2632
$ws.Range("I136").Value = 1543.7778
$ws.Range("J136").Value = 4250.364
$ws.Range("K136").Value = 4631.3334
$ws.Range("L136").Value = 12751.092
$ws.Range("M136").Value = -2081.3334
$ws.Range("N136").Value = -17851.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 928.7143
$ws.Range("J81").Value = 1016.6667
$ws.Range("L81").Value = 2033.3334
$ws.Range("N81").Value = -4155.3334

$ws.Range("H84").Value = 928.7143
$ws.Range("J84").Value = 1016.6667
$ws.Range("L84").Value = 10166.667
$ws.Range("N84").Value = -20774.667

$ws.Range("H107").Value = 537.7619
$ws.Range("I107").Value = 199
$ws.Range("J107").Value = 1088.25
$ws.Range("K107").Value = 597
$ws.Range("L107").Value = 3264.75
$ws.Range("M107").Value = 1323
$ws.Range("N107").Value = -7104.75

$ws.Range("H132").Value = 1312.0566
$ws.Range("I132").Value = 1014.75555
$ws.Range("J132").Value = 2984.375
$ws.Range("K132").Value = 3044.26665
$ws.Range("L132").Value = 8953.125
$ws.Range("M132").Value = -514.26665
$ws.Range("N132").Value = -14013.125

$ws.Range("H136").Value = 15433830
$ws.Range("I136").Value = 17362622
$ws.Range("K136").Value = 52087866
$ws.Range("M136").Value = -52085316
